$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Range("D130").Value = 44736
$ws.Range("H130").Value = "Madrigal"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 70
$ws.Range("K130").Value = 22000
$ws.Range("L130").Value = 23000
$ws.Range("M130").Value = 22429
$ws.Range("N130").Value = "`$/caja 30 unidades"
$ws.Range("O130").Value = "Provincia de Limarí"
$ws.Range("P130").Value = 748
$ws.Range("Q130").Value = 30

# Row 131
$ws.Range("D131").Value = 44390
$ws.Range("H131").Value = "Española"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 125
$ws.Range("K131").Value = 18000
$ws.Range("L131").Value = 18000
$ws.Range("M131").Value = 18000
$ws.Range("N131").Value = "`$/caja 30 unidades"
$ws.Range("O131").Value = "Provincia de Limarí"
$ws.Range("P131").Value = 600
$ws.Range("Q131").Value = 30

# Row 132
$ws.Range("D132").Value = 44384
$ws.Range("H132").Value = "Argentina(o)"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 40
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = 18000
$ws.Range("N132").Value = "`$/caja 50 unidades"
$ws.Range("O132").Value = "Provincia de Limarí"
$ws.Range("P132").Value = 360
$ws.Range("Q132").Value = 50

# Row 133
$ws.Range("D133").Value = 44384
$ws.Range("H133").Value = "Española"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 70
$ws.Range("K133").Value = 18000
$ws.Range("L133").Value = 19000
$ws.Range("M133").Value = 18429
$ws.Range("N133").Value = "`$/caja 30 unidades"
$ws.Range("O133").Value = "Provincia de Limarí"
$ws.Range("P133").Value = 614
$ws.Range("Q133").Value = 30

# Row 134
$ws.Range("D134").Value = 44411
$ws.Range("H134").Value = "Española"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 100
$ws.Range("K134").Value = 17000
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = 17500
$ws.Range("N134").Value = "`$/caja 30 unidades"
$ws.Range("O134").Value = "Provincia de Limarí"
$ws.Range("P134").Value = 583
$ws.Range("Q134").Value = 30

# Row 135
$ws.Range("D135").Value = 44411
$ws.Range("H135").Value = "Madrigal"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 80
$ws.Range("K135").Value = 17000
$ws.Range("L135").Value = 17000
$ws.Range("M135").Value = 17000
$ws.Range("N135").Value = "`$/caja 40 unidades"
$ws.Range("O135").Value = "Provincia de Limarí"
$ws.Range("P135").Value = 425
$ws.Range("Q135").Value = 40

# Row 136
$ws.Range("D136").Value = 44476
$ws.Range("H136").Value = "Española"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 100
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 12000
$ws.Range("N136").Value = "`$/caja 30 unidades"
$ws.Range("O136").Value = "Región Metropolitana"
$ws.Range("P136").Value = 400
$ws.Range("Q136").Value = 30

# Row 137
$ws.Range("D137").Value = 44476
$ws.Range("H137").Value = "Madrigal"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 80
$ws.Range("K137").Value = 12000
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = 12000
$ws.Range("N137").Value = "`$/caja 40 unidades"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 300
$ws.Range("Q137").Value = 40

# Row 138
$ws.Range("D138").Value = 44476
$ws.Range("H138").Value = "Madrigal"
$ws.Range("I138").Value = "Segunda"
$ws.Range("J138").Value = 30
$ws.Range("K138").Value = 10000
$ws.Range("L138").Value = 10000
$ws.Range("M138").Value = 10000
$ws.Range("N138").Value = "`$/caja 70 unidades"
$ws.Range("O138").Value = "Región Metropolitana"
$ws.Range("P138").Value = 143
$ws.Range("Q138").Value = 70

# Row 139
$ws.Range("D139").Value = 44482
$ws.Range("H139").Value = "Española"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 50
$ws.Range("K139").Value = 10000
$ws.Range("L139").Value = 10000
$ws.Range("M139").Value = 10000
$ws.Range("N139").Value = "`$/caja 30 unidades"
$ws.Range("O139").Value = "Región Metropolitana"
$ws.Range("P139").Value = 333
$ws.Range("Q139").Value = 30

# Row 140
$ws.Range("D140").Value = 44474
$ws.Range("H140").Value = "Española"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 50
$ws.Range("K140").Value = 13000
$ws.Range("L140").Value = 13000
$ws.Range("M140").Value = 13000
$ws.Range("N140").Value = "`$/caja 30 unidades"
$ws.Range("O140").Value = "Región Metropolitana"
$ws.Range("P140").Value = 433
$ws.Range("Q140").Value = 30

# Row 141
$ws.Range("D141").Value = 44474
$ws.Range("H141").Value = "Madrigal"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 50
$ws.Range("K141").Value = 12000
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 12000
$ws.Range("N141").Value = "`$/caja 40 unidades"
$ws.Range("O141").Value = "Región Metropolitana"
$ws.Range("P141").Value = 300
$ws.Range("Q141").Value = 40

# Row 142
$ws.Range("D142").Value = 44418
$ws.Range("H142").Value = "Madrigal"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 55
$ws.Range("K142").Value = 17000
$ws.Range("L142").Value = 17000
$ws.Range("M142").Value = 17000
$ws.Range("N142").Value = "`$/caja 40 unidades"
$ws.Range("O142").Value = "Provincia de Limarí"
$ws.Range("P142").Value = 425
$ws.Range("Q142").Value = 40

# Row 143
$ws.Range("D143").Value = 44494
$ws.Range("H143").Value = "Española"
$ws.Range("I143").Value = "Extra"
$ws.Range("J143").Value = 1000
$ws.Range("K143").Value = 500
$ws.Range("L143").Value = 500
$ws.Range("M143").Value = 500
$ws.Range("N143").Value = "`$/unidad"
$ws.Range("O143").Value = "Región Metropolitana"
$ws.Range("P143").Value = 500
$ws.Range("Q143").Value = 1

# Row 144
$ws.Range("D144").Value = 44398
$ws.Range("H144").Value = "Española"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 100
$ws.Range("K144").Value = 18000
$ws.Range("L144").Value = 18000
$ws.Range("M144").Value = 18000
$ws.Range("N144").Value = "`$/caja 30 unidades"
$ws.Range("O144").Value = "Provincia de Limarí"
$ws.Range("P144").Value = 600
$ws.Range("Q144").Value = 30

# Row 145
$ws.Range("D145").Value = 44398
$ws.Range("H145").Value = "Madrigal"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 40
$ws.Range("K145").Value = 18000
$ws.Range("L145").Value = 18000
$ws.Range("M145").Value = 18000
$ws.Range("N145").Value = "`$/caja 40 unidades"
$ws.Range("O145").Value = "Provincia de Limarí"
$ws.Range("P145").Value = 450
$ws.Range("Q145").Value = 40

# Row 146
$ws.Range("D146").Value = 44726
$ws.Range("H146").Value = "Madrigal"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 55
$ws.Range("K146").Value = 21000
$ws.Range("L146").Value = 21000
$ws.Range("M146").Value = 21000
$ws.Range("N146").Value = "`$/caja 30 unidades"
$ws.Range("O146").Value = "Provincia de Limarí"
$ws.Range("P146").Value = 700
$ws.Range("Q146").Value = 30

# Row 147
$ws.Range("D147").Value = 44426
$ws.Range("H147").Value = "Argentina(o)"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 80
$ws.Range("K147").Value = 16000
$ws.Range("L147").Value = 16000
$ws.Range("M147").Value = 16000
$ws.Range("N147").Value = "`$/caja 50 unidades"
$ws.Range("O147").Value = "Provincia de Limarí"
$ws.Range("P147").Value = 320
$ws.Range("Q147").Value = 50

# Row 148
$ws.Range("D148").Value = 44477
$ws.Range("H148").Value = "Madrigal"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 50
$ws.Range("K148").Value = 12000
$ws.Range("L148").Value = 12000
$ws.Range("M148").Value = 12000
$ws.Range("N148").Value = "`$/caja 40 unidades"
$ws.Range("O148").Value = "Región Metropolitana"
$ws.Range("P148").Value = 300
$ws.Range("Q148").Value = 40

# Row 149
$ws.Range("D149").Value = 44473
$ws.Range("H149").Value = "Española"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 80
$ws.Range("K149").Value = 12000
$ws.Range("L149").Value = 13000
$ws.Range("M149").Value = 12500
$ws.Range("N149").Value = "`$/caja 30 unidades"
$ws.Range("O149").Value = "Región Metropolitana"
$ws.Range("P149").Value = 417
$ws.Range("Q149").Value = 30

# Row 150
$ws.Range("D150").Value = 44473
$ws.Range("H150").Value = "Madrigal"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 100
$ws.Range("K150").Value = 12000
$ws.Range("L150").Value = 12000
$ws.Range("M150").Value = 12000
$ws.Range("N150").Value = "`$/caja 40 unidades"
$ws.Range("O150").Value = "Región Metropolitana"
$ws.Range("P150").Value = 300
$ws.Range("Q150").Value = 40

# Row 151
$ws.Range("D151").Value = 44326
$ws.Range("H151").Value = "Madrigal"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 80
$ws.Range("K151").Value = 18000
$ws.Range("L151").Value = 18000
$ws.Range("M151").Value = 18000
$ws.Range("N151").Value = "`$/caja 40 unidades"
$ws.Range("O151").Value = "Provincia de Limarí"
$ws.Range("P151").Value = 450
$ws.Range("Q151").Value = 40

# Row 152
$ws.Range("D152").Value = 44385
$ws.Range("H152").Value = "Argentina(o)"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 90
$ws.Range("K152").Value = 16500
$ws.Range("L152").Value = 17000
$ws.Range("M152").Value = 16778
$ws.Range("N152").Value = "`$/caja 50 unidades"
$ws.Range("O152").Value = "Provincia de Limarí"
$ws.Range("P152").Value = 336
$ws.Range("Q152").Value = 50

# Row 153
$ws.Range("D153").Value = 44385
$ws.Range("H153").Value = "Española"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 80
$ws.Range("K153").Value = 18000
$ws.Range("L153").Value = 18000
$ws.Range("M153").Value = 18000
$ws.Range("N153").Value = "`$/caja 30 unidades"
$ws.Range("O153").Value = "Provincia de Limarí"
$ws.Range("P153").Value = 600
$ws.Range("Q153").Value = 30

# Row 154
$ws.Range("D154").Value = 44385
$ws.Range("H154").Value = "Madrigal"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 80
$ws.Range("K154").Value = 17000
$ws.Range("L154").Value = 18000
$ws.Range("M154").Value = 17500
$ws.Range("N154").Value = "`$/caja 40 unidades"
$ws.Range("O154").Value = "Provincia de Limarí"
$ws.Range("P154").Value = 438
$ws.Range("Q154").Value = 40

# Row 155
$ws.Range("D155").Value = 44518
$ws.Range("H155").Value = "Madrigal"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 110
$ws.Range("K155").Value = 12000
$ws.Range("L155").Value = 12000
$ws.Range("M155").Value = 12000
$ws.Range("N155").Value = "`$/caja 40 unidades"
$ws.Range("O155").Value = "Región del Maule"
$ws.Range("P155").Value = 300
$ws.Range("Q155").Value = 40

# Row 156
$ws.Range("D156").Value = 44348
$ws.Range("H156").Value = "Española"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 75
$ws.Range("K156").Value = 18000
$ws.Range("L156").Value = 20000
$ws.Range("M156").Value = 19067
$ws.Range("N156").Value = "`$/caja 30 unidades"
$ws.Range("O156").Value = "Provincia de Limarí"
$ws.Range("P156").Value = 636
$ws.Range("Q156").Value = 30

# Row 157
$ws.Range("D157").Value = 44386
$ws.Range("H157").Value = "Argentina(o)"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 40
$ws.Range("K157").Value = 18000
$ws.Range("L157").Value = 18000
$ws.Range("M157").Value = 18000
$ws.Range("N157").Value = "`$/caja 50 unidades"
$ws.Range("O157").Value = "Provincia de Limarí"
$ws.Range("P157").Value = 360
$ws.Range("Q157").Value = 50

# Row 158
$ws.Range("D158").Value = 44386
$ws.Range("H158").Value = "Madrigal"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 70
$ws.Range("K158").Value = 17000
$ws.Range("L158").Value = 17000
$ws.Range("M158").Value = 17000
$ws.Range("N158").Value = "`$/caja 40 unidades"
$ws.Range("O158").Value = "Provincia de Limarí"
$ws.Range("P158").Value = 425
$ws.Range("Q158").Value = 40

# Row 159
$ws.Range("D159").Value = 44433
$ws.Range("H159").Value = "Española"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 110
$ws.Range("K159").Value = 15000
$ws.Range("L159").Value = 15000
$ws.Range("M159").Value = 15000
$ws.Range("N159").Value = "`$/caja 30 unidades"
$ws.Range("O159").Value = "Provincia de Limarí"
$ws.Range("P159").Value = 500
$ws.Range("Q159").Value = 30

# Row 160
$ws.Range("D160").Value = 44433
$ws.Range("H160").Value = "Española"
$ws.Range("I160").Value = "Segunda"
$ws.Range("J160").Value = 55
$ws.Range("K160").Value = 12000
$ws.Range("L160").Value = 12000
$ws.Range("M160").Value = 12000
$ws.Range("N160").Value = "`$/caja 30 unidades"
$ws.Range("O160").Value = "Provincia de Limarí"
$ws.Range("P160").Value = 400
$ws.Range("Q160").Value = 30

# Row 161
$ws.Range("D161").Value = 44421
$ws.Range("H161").Value = "Madrigal"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 95
$ws.Range("K161").Value = 17000
$ws.Range("L161").Value = 17000
$ws.Range("M161").Value = 17000
$ws.Range("N161").Value = "`$/caja 40 unidades"
$ws.Range("O161").Value = "Provincia de Limarí"
$ws.Range("P161").Value = 425
$ws.Range("Q161").Value = 40

# Row 162
$ws.Range("D162").Value = 44442
$ws.Range("H162").Value = "Argentina(o)"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 90
$ws.Range("K162").Value = 13000
$ws.Range("L162").Value = 13000
$ws.Range("M162").Value = 13000
$ws.Range("N162").Value = "`$/caja 50 unidades"
$ws.Range("O162").Value = "Región Metropolitana"
$ws.Range("P162").Value = 260
$ws.Range("Q162").Value = 50

# Row 163
$ws.Range("D163").Value = 44483
$ws.Range("H163").Value = "Española"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 150
$ws.Range("K163").Value = 14000
$ws.Range("L163").Value = 14000
$ws.Range("M163").Value = 14000
$ws.Range("N163").Value = "`$/caja 30 unidades"
$ws.Range("O163").Value = "Región Metropolitana"
$ws.Range("P163").Value = 467
$ws.Range("Q163").Value = 30

# Row 164
$ws.Range("D164").Value = 44187
$ws.Range("H164").Value = "Española"
$ws.Range("I164").Value = "Segunda"
$ws.Range("J164").Value = 40
$ws.Range("K164").Value = 15000
$ws.Range("L164").Value = 16000
$ws.Range("M164").Value = 15500
$ws.Range("N164").Value = "`$/caja 40 unidades"
$ws.Range("O164").Value = "Región del Maule"
$ws.Range("P164").Value = 388
$ws.Range("Q164").Value = 40

# New Row 165
$ws.Range("A165").Value = 10
$ws.Range("B165").Value = "Vega Modelo de Temuco"
$ws.Range("C165").Value = "La Araucanía"
$ws.Range("D165").Value = 44519
$ws.Range("D165").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E165").Value = 9
$ws.Range("F165").Value = 100112013
$ws.Range("G165").Value = "Alcachofa"
$ws.Range("H165").Value = "Madrigal"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 65
$ws.Range("K165").Value = 12000
$ws.Range("L165").Value = 12000
$ws.Range("M165").Value = 12000
$ws.Range("N165").Value = "`$/caja 40 unidades"
$ws.Range("O165").Value = "Región del Maule"
$ws.Range("P165").Value = 300
$ws.Range("Q165").Value = 40
$ws.Range("R165").Value = "Hortaliza"
